$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 390.80646
$ws.Range("I96").Value = 308.07144
$ws.Range("J96").Value = 1163
$ws.Range("K96").Value = 924.21432
$ws.Range("L96").Value = 3489
$ws.Range("M96").Value = 448.78568
$ws.Range("N96").Value = -6235
$ws.Range("H98").Value = 1147.2963
$ws.Range("I98").Value = 838.0909
$ws.Range("J98").Value = 2507.8
$ws.Range("K98").Value = 838.0909
$ws.Range("L98").Value = 2507.8
$ws.Range("M98").Value = 659.9091
$ws.Range("N98").Value = -5503.8
$ws.Range("H122").Value = 1147.2963
$ws.Range("I122").Value = 838.0909
$ws.Range("J122").Value = 2507.8
$ws.Range("K122").Value = 2514.2727
$ws.Range("L122").Value = 7523.400000000001
$ws.Range("M122").Value = -64.27269999999999
$ws.Range("N122").Value = -12423.4
$ws.Range("H125").Value = 5035
$ws.Range("I125").Value = 5034
$ws.Range("J125").Value = 5035.5
$ws.Range("K125").Value = 45306
$ws.Range("L125").Value = 45319.5
$ws.Range("M125").Value = -42846
$ws.Range("N125").Value = -50239.5
$ws.Range("H132").Value = 1709.05
$ws.Range("I132").Value = 1767.8125
$ws.Range("J132").Value = 1474
$ws.Range("K132").Value = 5303.4375
$ws.Range("L132").Value = 4422
$ws.Range("M132").Value = -2773.4375
$ws.Range("N132").Value = -9482

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2143.4565
$ws.Range("I2").Value = 1426.3125
$ws.Range("J2").Value = 3782.6428
$ws.Range("K2").Value = 1426.3125
$ws.Range("L2").Value = 3782.6428
$ws.Range("M2").Value = -1313.3125
$ws.Range("N2").Value = -4008.6428
$ws.Range("H45").Value = 8895.666999999999
$ws.Range("I45").Value = 10539.833
$ws.Range("K45").Value = 10539.833
$ws.Range("M45").Value = -10162.833
$ws.Range("H61").Value = 5237.7676
$ws.Range("I61").Value = 3125.639
$ws.Range("J61").Value = 16100.143
$ws.Range("K61").Value = 3125.639
$ws.Range("L61").Value = 16100.143
$ws.Range("M61").Value = -2913.639
$ws.Range("N61").Value = -16524.143
$ws.Range("H63").Value = 7802.737
$ws.Range("I63").Value = 2167.6453
$ws.Range("J63").Value = 32758.143
$ws.Range("K63").Value = 2167.6453
$ws.Range("L63").Value = 32758.143
$ws.Range("M63").Value = -1481.6453
$ws.Range("N63").Value = -34130.143
$ws.Range("H66").Value = 7802.737
$ws.Range("I66").Value = 2167.6453
$ws.Range("J66").Value = 32758.143
$ws.Range("K66").Value = 10838.2265
$ws.Range("L66").Value = 163790.715
$ws.Range("M66").Value = -7406.226500000001
$ws.Range("N66").Value = -170654.715
$ws.Range("H74").Value = 8431.766
$ws.Range("I74").Value = 9459.405000000001
$ws.Range("J74").Value = 4629.5
$ws.Range("K74").Value = 9459.405000000001
$ws.Range("L74").Value = 4629.5
$ws.Range("M74").Value = -8585.405000000001
$ws.Range("N74").Value = -6377.5
$ws.Range("H77").Value = 8431.766
$ws.Range("I77").Value = 9459.405000000001
$ws.Range("J77").Value = 4629.5
$ws.Range("K77").Value = 47297.025
$ws.Range("L77").Value = 23147.5
$ws.Range("M77").Value = -42929.025
$ws.Range("N77").Value = -31883.5
$ws.Range("H102").Value = 2412.1904
$ws.Range("I102").Value = 2412.1904
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2412.1904
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -790.1904
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 2143.4565
$ws.Range("I116").Value = 1426.3125
$ws.Range("J116").Value = 3782.6428
$ws.Range("K116").Value = 1426.3125
$ws.Range("L116").Value = 3782.6428
$ws.Range("M116").Value = 867.6875
$ws.Range("N116").Value = -8370.6428
$ws.Range("H132").Value = 3419.625
$ws.Range("I132").Value = 2945.3333
$ws.Range("K132").Value = 8835.999899999999
$ws.Range("M132").Value = -6305.999899999999
$ws.Range("H136").Value = 5237.7676
$ws.Range("I136").Value = 3125.639
$ws.Range("J136").Value = 16100.143
$ws.Range("K136").Value = 9376.917000000001
$ws.Range("L136").Value = 48300.429
$ws.Range("M136").Value = -6826.917000000001
$ws.Range("N136").Value = -53400.429

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2143.4565
$ws.Range("I3").Value = 1426.3125
$ws.Range("J3").Value = 3782.6428
$ws.Range("K3").Value = 1426.3125
$ws.Range("L3").Value = 3782.6428
$ws.Range("M3").Value = -1312.3125
$ws.Range("N3").Value = -4010.6428
$ws.Range("H80").Value = 759.625
$ws.Range("I80").Value = 160.33333
$ws.Range("K80").Value = 160.33333
$ws.Range("M80").Value = 837.6666700000001
$ws.Range("H83").Value = 759.625
$ws.Range("I83").Value = 160.33333
$ws.Range("K83").Value = 801.6666499999999
$ws.Range("M83").Value = 4190.33335
$ws.Range("H99").Value = 2764.5
$ws.Range("I99").Value = 2844.875
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 2844.875
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = -1346.875
$ws.Range("N99").Value = -4796
$ws.Range("H134").Value = 6145.707
$ws.Range("I134").Value = 2918.152
$ws.Range("J134").Value = 18518
$ws.Range("K134").Value = 8754.456
$ws.Range("L134").Value = 55554
$ws.Range("M134").Value = -6219.456
$ws.Range("N134").Value = -60624

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 249.44444
$ws.Range("I22").Value = 249.44444
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 249.44444
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 100.55556
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 2669.9355
$ws.Range("I58").Value = 1404.0385
$ws.Range("J58").Value = 9252.6
$ws.Range("K58").Value = 1404.0385
$ws.Range("L58").Value = 9252.6
$ws.Range("M58").Value = -1201.0385
$ws.Range("N58").Value = -9658.6
$ws.Range("H99").Value = 9037.591
$ws.Range("I99").Value = 4230.1665
$ws.Range("J99").Value = 12365.808
$ws.Range("K99").Value = 4230.1665
$ws.Range("L99").Value = 12365.808
$ws.Range("M99").Value = -2732.1665
$ws.Range("N99").Value = -15361.808
$ws.Range("H107").Value = 1150.5
$ws.Range("I107").Value = 1126.1666
$ws.Range("J107").Value = 1223.5
$ws.Range("K107").Value = 1126.1666
$ws.Range("L107").Value = 1223.5
$ws.Range("M107").Value = 793.8334
$ws.Range("N107").Value = -5063.5
$ws.Range("H122").Value = 1127.7391
$ws.Range("I122").Value = 1049.2941
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 3147.8823
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -697.8823000000002
$ws.Range("N122").Value = -8950
$ws.Range("H126").Value = 9037.591
$ws.Range("I126").Value = 4230.1665
$ws.Range("J126").Value = 12365.808
$ws.Range("K126").Value = 12690.4995
$ws.Range("L126").Value = 37097.424
$ws.Range("M126").Value = -10220.4995
$ws.Range("N126").Value = -42037.424
$ws.Range("H132").Value = 20395.014
$ws.Range("I132").Value = 11685.437
$ws.Range("J132").Value = 43205.81
$ws.Range("K132").Value = 35056.311
$ws.Range("L132").Value = 129617.43
$ws.Range("M132").Value = -32526.311
$ws.Range("N132").Value = -134677.43
$ws.Range("H134").Value = 3141.861
$ws.Range("I134").Value = 2541.0967
$ws.Range("J134").Value = 6866.6
$ws.Range("K134").Value = 7623.2901
$ws.Range("L134").Value = 20599.8
$ws.Range("M134").Value = -5088.2901
$ws.Range("N134").Value = -25669.8
$ws.Range("H136").Value = 2669.9355
$ws.Range("I136").Value = 1404.0385
$ws.Range("J136").Value = 9252.6
$ws.Range("K136").Value = 4212.1155
$ws.Range("L136").Value = 27757.8
$ws.Range("M136").Value = -1662.1155
$ws.Range("N136").Value = -32857.8

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6336.0557
$ws.Range("I70").Value = 3770.9167
$ws.Range("J70").Value = 11466.333
$ws.Range("K70").Value = 3770.9167
$ws.Range("L70").Value = 11466.333
$ws.Range("M70").Value = -3500.9167
$ws.Range("N70").Value = -12006.333
$ws.Range("H73").Value = 6336.0557
$ws.Range("I73").Value = 3770.9167
$ws.Range("J73").Value = 11466.333
$ws.Range("K73").Value = 3770.9167
$ws.Range("L73").Value = 11466.333
$ws.Range("M73").Value = -2834.9167
$ws.Range("N73").Value = -13338.333
$ws.Range("H107").Value = 287.09677
$ws.Range("I107").Value = 317.77274
$ws.Range("J107").Value = 212.11111
$ws.Range("K107").Value = 317.77274
$ws.Range("L107").Value = 212.11111
$ws.Range("M107").Value = 1602.22726
$ws.Range("N107").Value = -4052.11111
$ws.Range("H132").Value = 12794.212
$ws.Range("I132").Value = 9548.593000000001
$ws.Range("J132").Value = 27399.5
$ws.Range("K132").Value = 28645.779
$ws.Range("L132").Value = 82198.5
$ws.Range("M132").Value = -26115.779
$ws.Range("N132").Value = -87258.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1140.6666
$ws.Range("I22").Value = 875
$ws.Range("J22").Value = 1273.5
$ws.Range("K22").Value = 875
$ws.Range("L22").Value = 1273.5
$ws.Range("M22").Value = -580
$ws.Range("N22").Value = -1863.5
$ws.Range("H27").Value = 1140.6666
$ws.Range("I27").Value = 875
$ws.Range("J27").Value = 1273.5
$ws.Range("K27").Value = 875
$ws.Range("L27").Value = 1273.5
$ws.Range("M27").Value = -768
$ws.Range("N27").Value = -1487.5
$ws.Range("H93").Value = 2774.05
$ws.Range("I93").Value = 2893.8948
$ws.Range("J93").Value = 497
$ws.Range("K93").Value = 2893.8948
$ws.Range("L93").Value = 497
$ws.Range("M93").Value = -1645.8948
$ws.Range("N93").Value = -2993
$ws.Range("H132").Value = 3558.054
$ws.Range("I132").Value = 2741.4067
$ws.Range("K132").Value = 8224.2201
$ws.Range("M132").Value = -5694.2201
$ws.Range("H136").Value = 4167.575
$ws.Range("I136").Value = 3394.4546
$ws.Range("J136").Value = 7812.2856
$ws.Range("K136").Value = 10183.3638
$ws.Range("L136").Value = 23436.8568
$ws.Range("M136").Value = -7633.363799999999
$ws.Range("N136").Value = -28536.8568

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1100
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1100
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1100
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -3846
$ws.Range("H100").Value = 618.1177
$ws.Range("I100").Value = 638
$ws.Range("J100").Value = 300
$ws.Range("K100").Value = 1276
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -735
$ws.Range("N100").Value = -1682
$ws.Range("H109").Value = 79416.28999999999
$ws.Range("J109").Value = 79416.28999999999
$ws.Range("L109").Value = 79416.28999999999
$ws.Range("N109").Value = -82190.28999999999
$ws.Range("H126").Value = 4254.757
$ws.Range("I126").Value = 3368.2917
$ws.Range("J126").Value = 5891.3076
$ws.Range("K126").Value = 10104.8751
$ws.Range("L126").Value = 17673.9228
$ws.Range("M126").Value = -7634.875100000001
$ws.Range("N126").Value = -22613.9228
$ws.Range("H132").Value = 3220.43
$ws.Range("I132").Value = 3181.351
$ws.Range("J132").Value = 3832.6667
$ws.Range("K132").Value = 9544.053
$ws.Range("L132").Value = 11498.0001
$ws.Range("M132").Value = -7014.053
$ws.Range("N132").Value = -16558.0001
$ws.Range("H136").Value = 4446514.5
$ws.Range("I136").Value = 6061767
$ws.Range("J136").Value = 4570.9165
$ws.Range("K136").Value = 18185301
$ws.Range("L136").Value = 13712.7495
$ws.Range("M136").Value = -18182751
$ws.Range("N136").Value = -18812.7495
